$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 32
$ws.Range("F5").Value = 304
$ws.Range("F6").Value = 434
$ws.Range("F8").Value = 1966
$ws.Range("F10").Value = 30
$ws.Range("F12").Value = 1588
$ws.Range("F13").Value = 1588
$ws.Range("F14").Value = 1315
$ws.Range("F15").Value = 50
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 436
$ws.Range("F24").Value = 6957
$ws.Range("F25").Value = 7531
$ws.Range("F31").Value = 239
$ws.Range("F32").Value = 11
$ws.Range("F34").Value = 38
$ws.Range("F36").Value = 1372
$ws.Range("F40").Value = 673
$ws.Range("F44").Value = 207
$ws.Range("F47").Value = 120
$ws.Range("G8").Value = "不可售"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 22
$ws.Range("F5").Value = 49
$ws.Range("F8").Value = 5
$ws.Range("F17").Value = 278

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2541
$ws.Range("F4").Value = 249
$ws.Range("F5").Value = 111

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 111
$ws.Range("F6").Value = 22
$ws.Range("F7").Value = 434
$ws.Range("F9").Value = 1966
$ws.Range("F10").Value = 30
$ws.Range("F12").Value = 1588
$ws.Range("F13").Value = 1588
$ws.Range("F15").Value = 1315
$ws.Range("F16").Value = 50
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 436
$ws.Range("F23").Value = 6957
$ws.Range("F24").Value = 7531
$ws.Range("F27").Value = 239
$ws.Range("F28").Value = 11
$ws.Range("F32").Value = 5
$ws.Range("F36").Value = 673
$ws.Range("F43").Value = 207
$ws.Range("F46").Value = 120
$ws.Range("F49").Value = 278
$ws.Range("G9").Value = "不可售"
